$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 57.14035266666667
$ws.Range("H2").Value = 171.421058
$ws.Range("I2").Value = 0.7274038390747541
$ws.Range("J2").Value = 0.7274038390747541
$ws.Range("Q2").Value = 73.86740999168022
$ws.Range("R2").Value = 664.8066899251221
$ws.Range("S2").Value = 0.6564302542032383
$ws.Range("T2").Value = 0.6564302542032383

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 57.14035266666667
$ws.Range("H3").Value = 171.421058
$ws.Range("I3").Value = 0.7274038390747541
$ws.Range("J3").Value = 0.7274038390747541
$ws.Range("Q3").Value = 7.986583279356889
$ws.Range("R3").Value = 71.87924951421201
$ws.Range("S3").Value = 0.07097358487151587
$ws.Range("T3").Value = 0.07097358487151587

# Row 4
$ws.Range("I4").Value = 0.08622113322131104
$ws.Range("J4").Value = 0.08622113322131104
$ws.Range("S4").Value = 0.07780844334028869
$ws.Range("T4").Value = 0.07780844334028869

# Row 5
$ws.Range("I5").Value = 0.08622113322131104
$ws.Range("J5").Value = 0.08622113322131104
$ws.Range("S5").Value = 0.008412689881022352
$ws.Range("T5").Value = 0.008412689881022352

# Row 6
$ws.Range("I6").Value = 0.1863750277039348
$ws.Range("J6").Value = 0.1863750277039348
$ws.Range("S6").Value = 0.1681902132499696
$ws.Range("T6").Value = 0.1681902132499696

# Row 7
$ws.Range("I7").Value = 0.1863750277039348
$ws.Range("J7").Value = 0.1863750277039348
$ws.Range("S7").Value = 0.01818481445396516
$ws.Range("T7").Value = 0.01818481445396516
